$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Hunk 2: split " o excel" into " o " + "excel" runs, and drop the
# "_GoBack" bookmark that used to sit right after "excel" (it gets
# relocated by hunk 1 below). We delete a non-empty range that spans
# the bookmark so it is cleanly removed along with the old text.
# ------------------------------------------------------------------
$rngExcel = $d.Content
$found = $rngExcel.Find.Execute(" o excel para tratarlos con otras aplicaciones externas", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rngExcel.Delete()
    $rngExcel.InsertAfter(" o ")
    $rngExcel.Collapse(0)
    $rngExcel.InsertAfter("excel")
    $rngExcel.Collapse(0)
    $rngExcel.InsertAfter(" para tratarlos con otras aplicaciones externas")
}

# ------------------------------------------------------------------
# Hunk 1: replace the "Agrupación de usuarios..." bullet with the new
# "Automatizacion: ..." bullet, and re-create the "_GoBack" bookmark
# right after "La aplicación " (its new home in the updated text).
# ------------------------------------------------------------------
$para = $d.Paragraphs(5)
$newText = "Automatizacion: La aplicación permite automatizar el acceso a paginas web elegidas para extraer datos periódicamente en forma de XML u otros formatos a través de un crawler que implentamos."
$para.Range.Text = $newText

$rngBm = $d.Content
$found = $rngBm.Find.Execute("permite automatizar el acceso a", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $rngBm.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $rngBm) | Out-Null
}
